$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 238 (shifts the existing rows 238-269 down to 239-270,
# matching the dimension growing from A1:R269 to A1:R270).
$ws.Rows.Item(238).Insert()

# Populate the newly inserted row 238 with the new record's data.
# Columns that are constant across this whole data block (A,B,C,E,F,G,I,Q,R)
# are copied from the row directly above (row 237); the rest use the new
# values introduced by the edit.
$ws.Range("A238").Value = 5
$ws.Range("B238").Value = "Macroferia Regional de Talca"
$ws.Range("C238").Value = "Maule"
$ws.Range("D238").Value = 45180
$ws.Range("E238").Value = 7
$ws.Range("F238").Value = 100112031
$ws.Range("G238").Value = "Poroto verde"
$ws.Range("H238").Value = "Sin especificar"
$ws.Range("I238").Value = "Primera"
$ws.Range("J238").Value = 150
$ws.Range("K238").Value = 23000
$ws.Range("L238").Value = 23000
$ws.Range("M238").Value = 23000
$ws.Range("N238").Value = '$/malla 25 kilos'
$ws.Range("O238").Value = 'Perú'
$ws.Range("P238").Value = 920
$ws.Range("Q238").Value = 25
$ws.Range("R238").Value = "Hortaliza"
